$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.666.95"
$ws.Range("E2").Value = "'  +3.11%  "
$ws.Range("D3").Value = "'2.424.32"
$ws.Range("E3").Value = "'  +9.08%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'329.08"
$ws.Range("E5").Value = "'  +13.62%  "
$ws.Range("D6").Value = "'105.67"
$ws.Range("E6").Value = "'  -4.41%  "
$ws.Range("D7").Value = "'0.649"
$ws.Range("E7").Value = "'  +4.03%  "
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("D9").Value = "'0.676"
$ws.Range("E9").Value = "'  +13.00%  "
$ws.Range("D10").Value = "'42.64"
$ws.Range("E10").Value = "'  -1.94%  "
$ws.Range("E11").Value = "'  +4.55%  "
$ws.Range("D12").Value = "'8.72"
$ws.Range("E12").Value = "'  +1.46%  "
$ws.Range("E13").Value = "'  +4.01%  "
$ws.Range("D14").Value = "'17.40"
$ws.Range("E14").Value = "'  +17.10%  "
$ws.Range("E15").Value = "'  +3.05%  "
$ws.Range("D16").Value = "'2.791.60"
$ws.Range("E16").Value = "'  +9.04%  "
$ws.Range("D17").Value = "'2.430.87"
$ws.Range("E17").Value = "'  +9.28%  "
$ws.Range("D18").Value = "'43.716.32"
$ws.Range("E18").Value = "'  +3.12%  "
$ws.Range("B19").Value = "'Uniswap"
$ws.Range("C19").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'7.56"
$ws.Range("E19").Value = "'  +6.34%  "
$ws.Range("B20").Value = "'ShibaInu"
$ws.Range("C20").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000111"
$ws.Range("E20").Value = "'  +6.00%  "
$ws.Range("D21").Value = "'76.33"
$ws.Range("E21").Value = "'  +5.04%  "
$ws.Range("B22").Value = "'PancakeSwap"
$ws.Range("C22").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "'3.53"
$ws.Range("E22").Value = "'  +5.71%  "
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'274.86"
$ws.Range("E23").Value = "'  +19.56%  "
$ws.Range("D24").Value = "'2.48"
$ws.Range("E24").Value = "'  +3.08%  "
$ws.Range("D25").Value = "'9.67"
$ws.Range("E25").Value = "'  +8.21%  "
$ws.Range("D26").Value = "'12.13"
$ws.Range("E26").Value = "'  +6.52%  "
$ws.Range("E27").Value = "'  -0.09%  "
$ws.Range("D28").Value = "'3.98"
$ws.Range("E28").Value = "'  +0.37%  "
$ws.Range("D29").Value = "'23.15"
$ws.Range("E29").Value = "'  +11.31%  "
$ws.Range("B30").Value = "'Toncoin"
$ws.Range("C30").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "'  +1.37%  "
$ws.Range("B31").Value = "'Monero"
$ws.Range("C31").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'178.24"
$ws.Range("E31").Value = "'  +2.76%  "
$ws.Range("D32").Value = "'38.23"
$ws.Range("E32").Value = "'  +1.95%  "
$ws.Range("E33").Value = "'  +4.03%  "
$ws.Range("E34").Value = "'  +7.20%  "
$ws.Range("D35").Value = "'6.01"
$ws.Range("E35").Value = "'  +7.35%  "
$ws.Range("E36").Value = "'  +6.55%  "
$ws.Range("D37").Value = "'4.91"
$ws.Range("E37").Value = "'  -0.72%  "
$ws.Range("D38").Value = "'0.0374"
$ws.Range("E38").Value = "'  +0.76%  "
$ws.Range("D39").Value = "'4.08"
$ws.Range("E39").Value = "'  -3.17%  "
$ws.Range("E40").Value = "'  +5.40%  "
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "'  +21.76%  "
$ws.Range("E42").Value = "'  +25.86%  "
$ws.Range("D43").Value = "'128.81"
$ws.Range("E43").Value = "'  +27.19%  "
$ws.Range("D44").Value = "'0.238"
$ws.Range("E44").Value = "'  +3.23%  "
$ws.Range("D45").Value = "'70.47"
$ws.Range("E45").Value = "'  -4.30%  "
$ws.Range("D46").Value = "'12.85"
$ws.Range("E46").Value = "'  +4.96%  "
$ws.Range("E47").Value = "'  +0.09%  "
$ws.Range("D48").Value = "'9.80"
$ws.Range("E48").Value = "'  +16.79%  "
$ws.Range("D49").Value = "'5.75"
$ws.Range("E49").Value = "'  +7.95%  "
$ws.Range("D50").Value = "'87.97"
$ws.Range("E50").Value = "'  +70.43%  "
$ws.Range("D51").Value = "'1.34"
$ws.Range("E51").Value = "'  +5.38%  "
